$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 490.33334
$ws.Range("I32").Value = 481
$ws.Range("J32").Value = 495
$ws.Range("K32").Value = 481
$ws.Range("L32").Value = 495
$ws.Range("M32").Value = -155
$ws.Range("N32").Value = -1147
# Row 64
$ws.Range("H64").Value = 3087.875
$ws.Range("I64").Value = 2800
$ws.Range("J64").Value = 3375.75
$ws.Range("K64").Value = 2800
$ws.Range("L64").Value = 3375.75
$ws.Range("M64").Value = -2552
$ws.Range("N64").Value = -3871.75
# Row 67
$ws.Range("H67").Value = 3087.875
$ws.Range("I67").Value = 2800
$ws.Range("J67").Value = 3375.75
$ws.Range("K67").Value = 2800
$ws.Range("L67").Value = 3375.75
$ws.Range("M67").Value = -1942
$ws.Range("N67").Value = -5091.75
# Row 76
$ws.Range("H76").Value = 3217.923
$ws.Range("I76").Value = 3157.4243
$ws.Range("J76").Value = 3550.6667
$ws.Range("K76").Value = 3157.4243
$ws.Range("L76").Value = 3550.6667
$ws.Range("M76").Value = -2842.4243
$ws.Range("N76").Value = -4180.6667
# Row 79
$ws.Range("H79").Value = 3217.923
$ws.Range("I79").Value = 3157.4243
$ws.Range("J79").Value = 3550.6667
$ws.Range("K79").Value = 3157.4243
$ws.Range("L79").Value = 3550.6667
$ws.Range("M79").Value = -2065.4243
$ws.Range("N79").Value = -5734.6667
# Row 137
$ws.Range("H137").Value = 3766.9607
$ws.Range("I137").Value = 3235
$ws.Range("J137").Value = 5321.923
$ws.Range("K137").Value = 9705
$ws.Range("L137").Value = 15965.769
$ws.Range("M137").Value = -7155
$ws.Range("N137").Value = -21065.769

$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 300
$ws.Range("I5").Value = 100
$ws.Range("J5").Value = 500
$ws.Range("K5").Value = 100
$ws.Range("L5").Value = 500
$ws.Range("M5").Value = 12
$ws.Range("N5").Value = -724
# Row 32
$ws.Range("H32").Value = 18425.404
$ws.Range("I32").Value = 16307.789
$ws.Range("J32").Value = 19484.21
$ws.Range("K32").Value = 16307.789
$ws.Range("L32").Value = 19484.21
$ws.Range("M32").Value = -16020.789
$ws.Range("N32").Value = -20058.21
# Row 37
$ws.Range("H37").Value = 14575.917
$ws.Range("I37").Value = 5966.5557
$ws.Range("J37").Value = 40404
$ws.Range("K37").Value = 5966.5557
$ws.Range("L37").Value = 40404
$ws.Range("M37").Value = -5693.5557
$ws.Range("N37").Value = -40950
# Row 54
$ws.Range("H54").Value = 50000
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 50000
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 50000
$ws.Range("N54").Value = -51538
# Row 63
$ws.Range("H63").Value = 12595945
$ws.Range("I63").Value = 46172004
$ws.Range("J63").Value = 4922.5
$ws.Range("K63").Value = 46172004
$ws.Range("L63").Value = 4922.5
$ws.Range("M63").Value = -46171318
$ws.Range("N63").Value = -6294.5
# Row 66
$ws.Range("H66").Value = 12595945
$ws.Range("I66").Value = 46172004
$ws.Range("J66").Value = 4922.5
$ws.Range("K66").Value = 230860020
$ws.Range("L66").Value = 24612.5
$ws.Range("M66").Value = -230856588
$ws.Range("N66").Value = -31476.5
# Row 88
$ws.Range("H88").Value = 33334582
$ws.Range("I88").Value = 66666664
$ws.Range("J88").Value = 2500
$ws.Range("K88").Value = 66666664
$ws.Range("L88").Value = 2500
$ws.Range("M88").Value = -66666258
$ws.Range("N88").Value = -3312
# Row 91
$ws.Range("H91").Value = 33334582
$ws.Range("I91").Value = 66666664
$ws.Range("J91").Value = 2500
$ws.Range("K91").Value = 66666664
$ws.Range("L91").Value = 2500
$ws.Range("M91").Value = -66665260
$ws.Range("N91").Value = -5308
# Row 92
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
# Row 109
$ws.Range("H109").Value = 26050
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 26050
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 26050
$ws.Range("N109").Value = -28824

$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 300
$ws.Range("I4").Value = 100
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 100
$ws.Range("L4").Value = 500
$ws.Range("M4").Value = 15
$ws.Range("N4").Value = -730
# Row 15
$ws.Range("H15").Value = 25428.572
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 25428.572
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 25428.572
$ws.Range("N15").Value = -25882.572
# Row 30
$ws.Range("H30").Value = 21505
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 21505
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 21505
$ws.Range("N30").Value = -21755
$ws.Range("M30").ClearContents()
# Row 94
$ws.Range("H94").Value = 18520006
$ws.Range("I94").Value = 25001284
$ws.Range("J94").Value = 2070.1428
$ws.Range("K94").Value = 25001284
$ws.Range("L94").Value = 2070.1428
$ws.Range("M94").Value = -25000833
$ws.Range("N94").Value = -2972.1428
# Row 105
$ws.Range("H105").Value = 2570
$ws.Range("I105").Value = 2200
$ws.Range("J105").Value = 3125
$ws.Range("K105").Value = 2200
$ws.Range("L105").Value = 3125
$ws.Range("M105").Value = -453
$ws.Range("N105").Value = -6619

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4319.3877
$ws.Range("I31").Value = 1710.5333
$ws.Range("J31").Value = 5470.353
$ws.Range("K31").Value = 1710.5333
$ws.Range("L31").Value = 5470.353
$ws.Range("M31").Value = -1415.5333
$ws.Range("N31").Value = -6060.353
# Row 34
$ws.Range("H34").Value = 4319.3877
$ws.Range("I34").Value = 1710.5333
$ws.Range("J34").Value = 5470.353
$ws.Range("K34").Value = 1710.5333
$ws.Range("L34").Value = 5470.353
$ws.Range("M34").Value = -1508.5333
$ws.Range("N34").Value = -5874.353
# Row 62
$ws.Range("H62").Value = 45460324
$ws.Range("I62").Value = 250005000
$ws.Range("J62").Value = 5950.778
$ws.Range("K62").Value = 250005000
$ws.Range("L62").Value = 5950.778
$ws.Range("M62").Value = -250004376
$ws.Range("N62").Value = -7198.778
# Row 65
$ws.Range("H65").Value = 45460324
$ws.Range("I65").Value = 250005000
$ws.Range("J65").Value = 5950.778
$ws.Range("K65").Value = 1250025000
$ws.Range("L65").Value = 29753.89
$ws.Range("M65").Value = -1250021880
$ws.Range("N65").Value = -35993.89
# Row 99
$ws.Range("H99").Value = 5514.1665
$ws.Range("I99").Value = 2530
$ws.Range("J99").Value = 6508.8887
$ws.Range("K99").Value = 2530
$ws.Range("L99").Value = 6508.8887
$ws.Range("M99").Value = -1032
$ws.Range("N99").Value = -9504.8887
# Row 126
$ws.Range("H126").Value = 5514.1665
$ws.Range("I126").Value = 2530
$ws.Range("J126").Value = 6508.8887
$ws.Range("K126").Value = 7590
$ws.Range("L126").Value = 19526.6661
$ws.Range("M126").Value = -5120
$ws.Range("N126").Value = -24466.6661

$ws = $wb.Worksheets.Item("CUL")
# Row 34
$ws.Range("H34").Value = 14974.4
$ws.Range("I34").Value = 27747
$ws.Range("J34").Value = 10329.818
$ws.Range("K34").Value = 83241
$ws.Range("L34").Value = 30989.454
$ws.Range("M34").Value = -83157
$ws.Range("N34").Value = -31157.454
# Row 39
$ws.Range("H39").Value = 18080.889
$ws.Range("I39").Value = 16000
$ws.Range("J39").Value = 18203.295
$ws.Range("K39").Value = 48000
$ws.Range("L39").Value = 54609.88499999999
$ws.Range("M39").Value = -47706
$ws.Range("N39").Value = -55197.88499999999
# Row 80
$ws.Range("H80").Value = 3657.65
$ws.Range("I80").Value = 2000
$ws.Range("J80").Value = 3841.8333
$ws.Range("K80").Value = 6000
$ws.Range("L80").Value = 11525.4999
$ws.Range("M80").Value = -5064
$ws.Range("N80").Value = -13397.4999
# Row 83
$ws.Range("H83").Value = 3657.65
$ws.Range("I83").Value = 2000
$ws.Range("J83").Value = 3841.8333
$ws.Range("K83").Value = 18000
$ws.Range("L83").Value = 34576.4997
$ws.Range("M83").Value = -13320
$ws.Range("N83").Value = -43936.4997
# Row 94
$ws.Range("H94").Value = 2527.7778
$ws.Range("I94").Value = 1000
$ws.Range("J94").Value = 2718.75
$ws.Range("K94").Value = 3000
$ws.Range("L94").Value = 8156.25
$ws.Range("M94").Value = -2324
$ws.Range("N94").Value = -9508.25
# Row 100
$ws.Range("H100").Value = 4480.1665
$ws.Range("I100").Value = 125
$ws.Range("J100").Value = 4876.091
$ws.Range("K100").Value = 375
$ws.Range("L100").Value = 14628.273
$ws.Range("M100").Value = 436
$ws.Range("N100").Value = -16250.273
# Row 106
$ws.Range("H106").Value = 3000
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 3000
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 9000
$ws.Range("N106").Value = -10892
# Row 122
$ws.Range("H122").Value = 2742.9321
$ws.Range("I122").Value = 843.6667
$ws.Range("J122").Value = 3084.8
$ws.Range("K122").Value = 7593.0003
$ws.Range("L122").Value = 27763.2
$ws.Range("M122").Value = -5143.0003
$ws.Range("N122").Value = -32663.2

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 6024.727
$ws.Range("I70").Value = 5430
$ws.Range("J70").Value = 7610.6665
$ws.Range("K70").Value = 5430
$ws.Range("L70").Value = 7610.6665
$ws.Range("M70").Value = -5160
$ws.Range("N70").Value = -8150.6665
# Row 73
$ws.Range("H73").Value = 6024.727
$ws.Range("I73").Value = 5430
$ws.Range("J73").Value = 7610.6665
$ws.Range("K73").Value = 5430
$ws.Range("L73").Value = 7610.6665
$ws.Range("M73").Value = -4494
$ws.Range("N73").Value = -9482.666499999999
# Row 102
$ws.Range("H102").Value = 2506.327
$ws.Range("I102").Value = 2200.558
$ws.Range("J102").Value = 3967.2222
$ws.Range("K102").Value = 2200.558
$ws.Range("L102").Value = 3967.2222
$ws.Range("M102").Value = -578.558
$ws.Range("N102").Value = -7211.2222
# Row 123
$ws.Range("H123").Value = 18900.062
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 18900.062
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 18900.062
$ws.Range("N123").Value = -23800.062

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 4799.857
$ws.Range("I7").Value = 2520
$ws.Range("J7").Value = 10499.5
$ws.Range("K7").Value = 2520
$ws.Range("L7").Value = 10499.5
$ws.Range("M7").Value = -2408
$ws.Range("N7").Value = -10723.5
# Row 46
$ws.Range("H46").Value = 2500.0667
$ws.Range("I46").Value = 2720.2
$ws.Range("J46").Value = 2390
$ws.Range("K46").Value = 2720.2
$ws.Range("L46").Value = 2390
$ws.Range("M46").Value = -2532.2
$ws.Range("N46").Value = -2766
# Row 126
$ws.Range("H126").Value = 4799.857
$ws.Range("I126").Value = 2520
$ws.Range("J126").Value = 10499.5
$ws.Range("K126").Value = 7560
$ws.Range("L126").Value = 31498.5
$ws.Range("M126").Value = -5090
$ws.Range("N126").Value = -36438.5

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 18908540
$ws.Range("I81").Value = 18908540
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 37817080
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -37816019
$ws.Range("N81").ClearContents()
# Row 84
$ws.Range("H84").Value = 18908540
$ws.Range("I84").Value = 18908540
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 189085400
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -189080096
$ws.Range("N84").ClearContents()
# Row 126
$ws.Range("H126").Value = 535184.3
$ws.Range("I126").Value = 2008.9
$ws.Range("J126").Value = 1068359.8
$ws.Range("K126").Value = 6026.700000000001
$ws.Range("L126").Value = 3205079.4
$ws.Range("M126").Value = -3556.700000000001
$ws.Range("N126").Value = -3210019.4
